$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.586.01"
$ws.Range("E2").Value = "  +5.88%  "
$ws.Range("D3").Value = "'1.918.70"
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").Value = "'335.21"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("D7").Value = "'0.4673"
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("E8").Value = "  +6.60%  "
$ws.Range("D9").Value = "'48.20"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "'0.08041"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").Value = "'22.40"
$ws.Range("E12").Value = "  +6.13%  "
$ws.Range("D13").Value = "'1.946.39"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").Value = "'6.011"
$ws.Range("E14").Value = "  +4.58%  "
$ws.Range("D15").Value = "'7.194"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").Value = "'89.89"
$ws.Range("E16").Value = "  +3.64%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "'0.06599"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").Value = "'17.87"
$ws.Range("E20").Value = "  +6.17%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").Value = "'29.557.77"
$ws.Range("E22").Value = "  +5.68%  "
$ws.Range("E23").Value = "  +5.63%  "
$ws.Range("D24").Value = "'11.65"
$ws.Range("E24").Value = "  +10.57%  "
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "'2.133.44"
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").Value = "'156.55"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +4.64%  "
$ws.Range("D29").Value = "'2.146"
$ws.Range("E29").Value = "  +6.36%  "
$ws.Range("D30").Value = "'5.743"
$ws.Range("E30").Value = "  +10.50%  "
$ws.Range("D31").Value = "'117.51"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "'1.078"
$ws.Range("E32").Value = "  +16.53%  "
$ws.Range("D33").Value = "'0.09485"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").Value = "'1.433"
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("E35").Value = "  +5.68%  "
$ws.Range("D36").Value = "'3.531"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").Value = "'0.06149"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").Value = "'0.02273"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("D39").Value = "'8.441"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("D40").Value = "'1.180"
$ws.Range("E40").Value = "  +3.24%  "
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").Value = "'0.1847"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("D43").Value = "'10.22"
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.341"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.243"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "'0.07512"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("E47").Value = "  +4.87%  "
$ws.Range("D48").Value = "'12.24"
$ws.Range("E48").Value = "  +4.37%  "
$ws.Range("D49").Value = "'1.937"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("D50").Value = "'113.31"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").Value = "'0.2991"
$ws.Range("E51").Value = "  +15.08%  "

Write-Host "Applied cryptos list update"
